# Auto-generated edit script applying scheduled market-data refresh to Zodiark_Profits sheets
# (workbook sheet tabs: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 45
$ws.Range("H45").Value = 2000
$ws.Range("J45").Value = 2000
$ws.Range("L45").Value = 6000
$ws.Range("N45").Value = -6384

# Row 47
$ws.Range("H47").Value = 13276.5
$ws.Range("I47").Value = 500
$ws.Range("K47").Value = 500
$ws.Range("M47").Value = 472

# Row 51
$ws.Range("H51").Value = 4778
$ws.Range("I51").Value = 4296
$ws.Range("J51").Value = 4821.8184
$ws.Range("K51").Value = 4296
$ws.Range("L51").Value = 4821.8184
$ws.Range("M51").Value = -3812
$ws.Range("N51").Value = -5789.8184

# Row 52
$ws.Range("H52").Value = 1385.2727
$ws.Range("I52").Value = 523.8
$ws.Range("J52").Value = 10000
$ws.Range("K52").Value = 1571.4
$ws.Range("L52").Value = 30000
$ws.Range("M52").Value = -1411.4
$ws.Range("N52").Value = -30320

# Row 59
$ws.Range("H59").Value = 10000
$ws.Range("J59").Value = 10000
$ws.Range("L59").Value = 30000
$ws.Range("N59").Value = -31114

# Row 86
$ws.Range("H86").Value = 1420.9445
$ws.Range("J86").Value = 1929.7142
$ws.Range("L86").Value = 1929.7142
$ws.Range("N86").Value = -4175.7142

# Row 89
$ws.Range("H89").Value = 1420.9445
$ws.Range("J89").Value = 1929.7142
$ws.Range("L89").Value = 9648.571
$ws.Range("N89").Value = -20880.571

# Row 103
$ws.Range("H103").Value = 284.2
$ws.Range("I103").Value = 348.14285
$ws.Range("J103").Value = 249.76923
$ws.Range("K103").Value = 1044.42855
$ws.Range("L103").Value = 749.30769
$ws.Range("M103").Value = -458.4285500000001
$ws.Range("N103").Value = -1921.30769

# Row 112
$ws.Range("H112").Value = 1533.3
$ws.Range("I112").Value = 414.5
$ws.Range("J112").Value = 1705.4231
$ws.Range("K112").Value = 1243.5
$ws.Range("L112").Value = 5116.2693
$ws.Range("M112").Value = -135.5
$ws.Range("N112").Value = -7332.2693

# Row 137
$ws.Range("H137").Value = 2077.92
$ws.Range("I137").Value = 2293
$ws.Range("J137").Value = 1620.875
$ws.Range("K137").Value = 6879
$ws.Range("L137").Value = 4862.625
$ws.Range("M137").Value = -4329
$ws.Range("N137").Value = -9962.625

# Row 141
$ws.Range("H141").Value = 3185.963
$ws.Range("I141").Value = 3234.6924
$ws.Range("J141").Value = 1919
$ws.Range("K141").Value = 9704.0772
$ws.Range("L141").Value = 5757
$ws.Range("M141").Value = -4524.0772
$ws.Range("N141").Value = -16117


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2483.2808
$ws.Range("I32").Value = 1195.3062
$ws.Range("K32").Value = 1195.3062
$ws.Range("M32").Value = -908.3062

# Row 45
$ws.Range("H45").Value = 1847.1072
$ws.Range("I45").Value = 1684.579
$ws.Range("J45").Value = 2190.2222
$ws.Range("K45").Value = 1684.579
$ws.Range("L45").Value = 2190.2222
$ws.Range("M45").Value = -1307.579
$ws.Range("N45").Value = -2944.2222

# Row 97
$ws.Range("H97").Value = 565.05884
$ws.Range("I97").Value = 390
$ws.Range("J97").Value = 720.6667
$ws.Range("K97").Value = 390
$ws.Range("L97").Value = 720.6667
$ws.Range("M97").Value = 106
$ws.Range("N97").Value = -1712.6667

# Row 122
$ws.Range("H122").Value = 4372.645
$ws.Range("I122").Value = 4526.85
$ws.Range("J122").Value = 4092.2727
$ws.Range("K122").Value = 13580.55
$ws.Range("L122").Value = 12276.8181
$ws.Range("M122").Value = -11130.55
$ws.Range("N122").Value = -17176.8181


$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 3382
$ws.Range("I107").Value = 2151.9473
$ws.Range("J107").Value = 6720.7144
$ws.Range("K107").Value = 2151.9473
$ws.Range("L107").Value = 6720.7144
$ws.Range("M107").Value = -231.9472999999998
$ws.Range("N107").Value = -10560.7144


$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 4286.8647
$ws.Range("J99").Value = 5048.6
$ws.Range("L99").Value = 5048.6
$ws.Range("N99").Value = -8044.6

# Row 126
$ws.Range("H126").Value = 4286.8647
$ws.Range("J126").Value = 5048.6
$ws.Range("L126").Value = 15145.8
$ws.Range("N126").Value = -20085.8


$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 447.75
$ws.Range("I33").Value = 347
$ws.Range("J33").Value = 467.9
$ws.Range("K33").Value = 2082
$ws.Range("L33").Value = 2807.4
$ws.Range("M33").Value = -1799
$ws.Range("N33").Value = -3373.4

# Row 35
$ws.Range("H35").Value = 600
$ws.Range("J35").Value = 600
$ws.Range("L35").Value = 1800
$ws.Range("N35").Value = -2376

# Row 41
$ws.Range("H41").Value = 2999
$ws.Range("I41").Value = 2999
$ws.Range("J41").Value = 2999
$ws.Range("K41").Value = 8997
$ws.Range("L41").Value = 8997
$ws.Range("M41").Value = -8659
$ws.Range("N41").Value = -9673

# Row 54
$ws.Range("H54").Value = 5074.25
$ws.Range("J54").Value = 5432.3335
$ws.Range("L54").Value = 16297.0005
$ws.Range("N54").Value = -17415.0005

# Row 59
$ws.Range("H59").Value = 9585.571
$ws.Range("I59").Value = 3749.5
$ws.Range("J59").Value = 11920
$ws.Range("K59").Value = 11248.5
$ws.Range("L59").Value = 35760
$ws.Range("M59").Value = -10708.5
$ws.Range("N59").Value = -36840

# Row 62
$ws.Range("H62").Value = 1125
$ws.Range("I62").Value = 1000
$ws.Range("J62").Value = 1500
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 4500
$ws.Range("M62").Value = -2314
$ws.Range("N62").Value = -5872

# Row 65
$ws.Range("H65").Value = 1125
$ws.Range("I65").Value = 1000
$ws.Range("J65").Value = 1500
$ws.Range("K65").Value = 9000
$ws.Range("L65").Value = 13500
$ws.Range("M65").Value = -5568
$ws.Range("N65").Value = -20364

# Row 69
$ws.Range("H69").Value = 6587.5
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

# Row 72
$ws.Range("H72").Value = 6587.5
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

# Row 120
$ws.Range("H120").Value = 10714.7
$ws.Range("I120").Value = 5429.4
$ws.Range("K120").Value = 16288.2
$ws.Range("M120").Value = -11450.2

# Row 140
$ws.Range("H140").Value = 2681
$ws.Range("I140").Value = 1944.2858
$ws.Range("K140").Value = 5832.857400000001
$ws.Range("M140").Value = -652.8574000000008


$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 4067.7144
$ws.Range("I102").Value = 3815.4
$ws.Range("J102").Value = 4698.5
$ws.Range("K102").Value = 3815.4
$ws.Range("L102").Value = 4698.5
$ws.Range("M102").Value = -2193.4
$ws.Range("N102").Value = -7942.5

# Row 122
$ws.Range("H122").Value = 3971.1538
$ws.Range("I122").Value = 2227.5
$ws.Range("J122").Value = 6761
$ws.Range("K122").Value = 6682.5
$ws.Range("L122").Value = 20283
$ws.Range("M122").Value = -4232.5
$ws.Range("N122").Value = -25183

# Row 126
$ws.Range("H126").Value = 5689.778
$ws.Range("J126").Value = 4458.2856
$ws.Range("L126").Value = 13374.8568
$ws.Range("N126").Value = -18314.8568


$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 1884.7333
$ws.Range("I113").Value = 1452.3636
$ws.Range("J113").Value = 3073.75
$ws.Range("K113").Value = 4357.0908
$ws.Range("L113").Value = 9221.25
$ws.Range("M113").Value = -2187.0908
$ws.Range("N113").Value = -13561.25

